$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate header row (row 1) ---
$ws.Range("A1").Value = "NOME"
$ws.Range("B1").Value = "IDADE"
$ws.Range("C1").Value = "ESTADO CÍVIL"
$ws.Range("D1").Value = "CIDADE"

# --- Populate data row (row 2) ---
$ws.Range("A2").Value = "Matheus"
$ws.Range("B2").Value = "17"
$ws.Range("C2").Value = "Solteiro"
$ws.Range("D2").Value = "São Paulo"

# --- Turn the range into an Excel Table (ListObject), using the existing header row ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D2"), $null, 1)
$lo.Name = "Tabela1"
$lo.TableStyle = "TableStyleLight9"

# --- Column widths for the new columns ---
$ws.Columns("C").ColumnWidth = 13.25
$ws.Columns("D").ColumnWidth = 8.25

# --- Selection / active cell as left by the editing session ---
$excel.Goto($ws.Range("C3:E6"))
